# Updates Price (D) and Volume(1h) (E) columns on Sheet1 to match the
# latest cryptos snapshot. D-column writes go through a NumberFormat="@"
# round-trip so numeric-looking text (e.g. "582.89") is not silently
# coerced to a Number by Excel; Style is then reset to "Normal" so no
# spurious formatting diff is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.119.08'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.467.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.87%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("E9").Value = '  +1.84%  '
$ws.Range("E10").Value = '  +0.44%  '
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.333'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.37'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.002.43'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.458.20'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.24%  '
$ws.Range("E19").Value = '  -2.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '348.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '3.99'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.31%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.593.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0901'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '500.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.68'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.10'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.61%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.327'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.82'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.39'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '142.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("E50").Value = '  -1.55%  '
$ws.Range("E51").Value = '  +0.02%  '
